$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.244.65'
$ws.Range("E2").Value = '  -1.60%  '
$ws.Range("D3").Value = '2.248.25'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''247.04'
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("D6").Value = '''0.626'
$ws.Range("E6").Value = '  -1.38%  '
$ws.Range("D7").Value = '''74.57'
$ws.Range("E7").Value = '  -5.61%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.620'
$ws.Range("E9").Value = '  -4.11%  '
$ws.Range("D10").Value = '''42.27'
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("E11").Value = '  -2.87%  '
$ws.Range("E12").Value = '  -2.20%  '
$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").Value = '2.240.01'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = '42.170.61'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '''72.48'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '''6.13'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").Value = '''2.22'
$ws.Range("E21").Value = '  +2.40%  '
$ws.Range("D22").Value = '''231.82'
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").Value = '''8.88'
$ws.Range("E23").Value = '  +37.44%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '''11.44'
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").Value = '''3.61'
$ws.Range("E26").Value = '  -4.97%  '
$ws.Range("D27").Value = '''2.30'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = '''2.17'
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").Value = '''169.25'
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = '''0.0828'
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").Value = '''0.124'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").Value = '''31.23'
$ws.Range("E33").Value = '  +3.41%  '
$ws.Range("E34").Value = '  -1.78%  '
$ws.Range("D35").Value = '''5.26'
$ws.Range("E35").Value = '  +10.08%  '
$ws.Range("E36").Value = '  -1.62%  '
$ws.Range("E37").Value = '  +3.73%  '
$ws.Range("D38").Value = '''14.05'
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("E39").Value = '  -3.47%  '
$ws.Range("D40").Value = '''5.78'
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("D42").Value = '''62.19'
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").Value = '''107.12'
$ws.Range("E43").Value = '  -5.16%  '
$ws.Range("E44").Value = '  +1.56%  '
$ws.Range("D45").Value = '''8.70'
$ws.Range("E45").Value = '  -2.37%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("E47").Value = '  -2.56%  '
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").Value = '''4.16'
$ws.Range("E50").Value = '  -9.90%  '
$ws.Range("E51").Value = '  -3.08%  '
